$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 134, pushing the existing rows 134-174 down to
# 136-176 (this reproduces the "shift down by 2" pattern seen across the
# whole diff automatically, without having to rewrite every row by hand).
$ws.Rows("134:135").Insert()

# --- Row 134: new weekly entry (Primera) ---
$ws.Cells.Item(134, 1).Value = 10
$ws.Cells.Item(134, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(134, 3).Value = "La Araucanía"
$ws.Cells.Item(134, 4).Value = 44551
$ws.Cells.Item(134, 5).Value = 9
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100102
$ws.Cells.Item(134, 8).Value = "Cítricos"
$ws.Cells.Item(134, 9).Value = 100102006
$ws.Cells.Item(134, 10).Value = "Pomelo"
$ws.Cells.Item(134, 11).Value = "Start Ruby"
$ws.Cells.Item(134, 12).Value = "Primera"
$ws.Cells.Item(134, 13).Value = 50
$ws.Cells.Item(134, 14).Value = 14000
$ws.Cells.Item(134, 15).Value = 14000
$ws.Cells.Item(134, 16).Value = 14000
$ws.Cells.Item(134, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(134, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(134, 19).Value = 933
$ws.Cells.Item(134, 20).Value = 15

# --- Row 135: new weekly entry (Segunda) ---
$ws.Cells.Item(135, 1).Value = 10
$ws.Cells.Item(135, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(135, 3).Value = "La Araucanía"
$ws.Cells.Item(135, 4).Value = 44551
$ws.Cells.Item(135, 5).Value = 9
$ws.Cells.Item(135, 6).Value = "Fruta"
$ws.Cells.Item(135, 7).Value = 100102
$ws.Cells.Item(135, 8).Value = "Cítricos"
$ws.Cells.Item(135, 9).Value = 100102006
$ws.Cells.Item(135, 10).Value = "Pomelo"
$ws.Cells.Item(135, 11).Value = "Start Ruby"
$ws.Cells.Item(135, 12).Value = "Segunda"
$ws.Cells.Item(135, 13).Value = 45
$ws.Cells.Item(135, 14).Value = 10000
$ws.Cells.Item(135, 15).Value = 10000
$ws.Cells.Item(135, 16).Value = 10000
$ws.Cells.Item(135, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(135, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(135, 19).Value = 667
$ws.Cells.Item(135, 20).Value = 15
